# "build off and def matrix"
# Fill in the scouting / defense_capability cells of the off/def matrix
# (columns D:F, i.e. unit_1r/unit_2r/unit_3r) on the "blue" sheet that were
# previously left at 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("blue")

# Block 1 (unit_1b) - scouting row
$ws.Range("D13:F13").Value = 1

# Block 2 (unit_2b) - defense_capability row
$ws.Range("D26:F26").Value = 1

# Block 2 (unit_2b) - scouting row
$ws.Range("D28:F28").Value = 1

# Block 3 (unit_3b) - scouting row
$ws.Range("D43:F43").Value = 1

# Restore the selection left on the sheet after the edit.
$ws.Range("F28").Select()
